{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"unchanging bits being black and ground brown\") !== -1) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Anchor paragraph not found\");\n}\n\nconst p1 = anchor.insertParagraph(\n  \"Idea: use 16*16 grid based search to simplify block search (won\\u2019t work for Goombas and the like)\",\n  Word.InsertLocation.after\n);\n\nconst p2 = p1.insertParagraph(\n  \"Idea: categorise objects as just \\u2018enemy\\u2019, \\u2018hazard\\u2019, \\u2018block\\u2019 etc.\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph that ends the ground-theme colour description; the\n# two new bullet items belong right after it.\n$searchText = \"unchanging bits being black and ground brown\"\n$findRng = $d.Content\n$found = $findRng.Find.Execute($searchText)\nif (-not $found) {\n    throw \"Anchor text not found\"\n}\n\n# Resolve the Find hit to its 1-based paragraph index in the document so we\n# can address the (soon to be inserted) sibling paragraphs without relying\n# on brittle hard-coded numbers.\n$anchorIndex = $d.Range(0, $findRng.Start).Paragraphs.Count\n$anchorPara = $d.Paragraphs.Item($anchorIndex)\n\n# Insert a new paragraph right after the anchor; Word carries the\n# ListParagraph style / numbering from the preceding paragraph onto it\n# automatically, same as the existing bullet items.\n$insertPoint = $anchorPara.Range\n$insertPoint.Collapse(0)\n$insertPoint.InsertParagraphAfter()\n\n$rightSingleQuote = [char]0x2019\n$leftSingleQuote = [char]0x2018\n\n$p1 = $d.Paragraphs.Item($anchorIndex + 1)\n$p1.Range.Text = \"Idea: use 16*16 grid based search to simplify block search (won\" + $rightSingleQuote + \"t work for Goombas and the like)\"\n\n$p1.Range.InsertParagraphAfter()\n\n$p2 = $d.Paragraphs.Item($anchorIndex + 2)\n$p2.Range.Text = \"Idea: categorise objects as just \" + $leftSingleQuote + \"enemy\" + $rightSingleQuote + \", \" + $leftSingleQuote + \"hazard\" + $rightSingleQuote + \", \" + $leftSingleQuote + \"block\" + $rightSingleQuote + \" etc.\"\n"}
